# Weekly refresh: a new week of "Perejil" price data is inserted as the
# new first data row right before the old top entry (row 29 in the
# "general" header-bearing sheet), pushing every following row down by one
# and appending the previously-last row as a brand new final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 29 - this shifts rows 29..140 down to 30..141,
# automatically carrying their existing Fecha/Volumen/Precio values with
# them (which is exactly the "shift by one" pattern seen across the diff).
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with this week's record.
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = "Terminal La Palmera de La Serena"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44677
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112044
$ws.Range("G29").Value = "Perejil"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 2400
$ws.Range("K29").Value = 2500
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = 2750
$ws.Range("N29").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O29").Value = "Provincia del Elquí"
$ws.Range("P29").Value = 1833
$ws.Range("Q29").Value = 1.5
$ws.Range("R29").Value = "Hortaliza"
